$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping")

# Add the new "region" column header
$ws.Range("D1").Value = "region"

# Classify each county as urban or rural.
# (Assigning all "urban" rows before any "rural" row so the
#  shared-string table is built up in the same order as the source file.)
$urbanRows = @(10,14,19,22,26,32,44,48,49,51,58,77,78,79,84)
foreach ($r in $urbanRows) {
    $ws.Range("D$r").Value = "urban"
}

$ruralRows = @(2,3,4,5,6,7,8,9,11,12,13,15,16,17,18,20,21,23,24,25,27,28,29,30,31,33,34,35,36,37,38,39,40,41,42,43,45,46,47,50,52,53,54,55,56,57,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,80,81,82,83,85,86,87,88,89)
foreach ($r in $ruralRows) {
    $ws.Range("D$r").Value = "rural"
}

# Reset the view: scroll back to the top and select the new header cell
$ws.Range("D1").Select()

# Re-create the hidden _FilterDatabase defined name scoped to the "mapping"
# sheet, covering the full used range including the new column
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=mapping!`$A`$1:`$D`$89")
$filterName.Visible = $false

